$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'28.162.77"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').Value = "'1.868.91"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.05%  '
$ws.Range('D5').Value = "'311.63"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('D6').Value = "'1.001"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').Value = "'0.5049"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.34%  '
$ws.Range('D8').Value = "'0.3920"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.32%  '
$ws.Range('D9').Value = "'0.09647"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.29%  '
$ws.Range('D10').Value = "'1.138"
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Value = "'40.89"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.28%  '
$ws.Range('D12').Value = "'6.497"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.15%  '
$ws.Range('D13').Value = "'20.93"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.16%  '
$ws.Range('D14').Value = "'1.867.07"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.64%  '
$ws.Range('D15').Value = "'1.002"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.25%  '
$ws.Range('E16').Value = '  -0.01%  '
$ws.Range('D17').Value = "'0.00001128"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.98%  '
$ws.Range('D18').Value = "'92.96"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.58%  '
$ws.Range('D19').Value = "'0.06621"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('D20').Value = "'17.52"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.75%  '
$ws.Range('E21').Value = '  +0.15%  '
$ws.Range('D22').Value = "'6.162"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.64%  '
$ws.Range('D23').Value = "'28.226.44"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('E24').Value = '  +1.29%  '
$ws.Range('D25').Value = "'2.279"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.38%  '
$ws.Range('D26').Value = "'2.536"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.34%  '
$ws.Range('D27').Value = "'2.089.04"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.63%  '
$ws.Range('D28').Value = "'21.20"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.51%  '
$ws.Range('D29').Value = "'157.90"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.25%  '
$ws.Range('D30').Value = "'127.39"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.96%  '
$ws.Range('D31').Value = "'0.1058"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.67%  '
$ws.Range('D32').Value = "'1.066"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.10%  '
$ws.Range('D33').Value = "'5.626"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.50%  '
$ws.Range('D34').Value = "'3.625"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.42%  '
$ws.Range('D35').Value = "'9.577"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.06%  '
$ws.Range('D36').Value = "'0.06735"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.45%  '
$ws.Range('D37').Value = "'0.02383"
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Value = "'0.2176"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.05%  '
$ws.Range('E39').Value = '  -1.37%  '
$ws.Range('D40').Value = "'0.6348"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.68%  '
$ws.Range('D41').Value = "'4.976"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.23%  '
$ws.Range('D42').Value = "'1.176"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.65%  '
$ws.Range('E43').Value = '  +0.18%  '
$ws.Range('D44').Value = "'13.54"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.78%  '
$ws.Range('D45').Value = "'0.6013"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('D46').Value = "'3.660"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.41%  '
$ws.Range('D47').Value = "'1.262"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.09%  '
$ws.Range('D48').Value = "'124.29"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.31%  '
$ws.Range('D49').Value = "'1.993"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.54%  '
$ws.Range('D50').Value = "'1.194"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.13%  '
$ws.Range('E51').Value = '  +0.71%  '
